$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 11323
$ws1.Range("F4").Value = 1759
$ws1.Range("F5").Value = 489
$ws1.Range("F6").Value = 781
$ws1.Range("F7").Value = 2362
$ws1.Range("F8").Value = 712
$ws1.Range("F9").Value = 945
$ws1.Range("F10").Value = 546
$ws1.Range("F11").Value = 409
$ws1.Range("F12").Value = 454
$ws1.Range("F13").Value = 418
$ws1.Range("F15").Value = 581
$ws1.Range("F16").Value = 41
$ws1.Range("F17").Value = 934
$ws1.Range("F18").Value = 406
$ws1.Range("F19").Value = 599
$ws1.Range("F20").Value = 949
$ws1.Range("F21").Value = 180
$ws1.Range("F22").Value = 893
$ws1.Range("F24").Value = 181
$ws1.Range("F27").Value = 627
$ws1.Range("F28").Value = 128
$ws1.Range("F29").Value = 71
$ws1.Range("F30").Value = 294

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 705
$ws2.Range("F4").Value = 84
$ws2.Range("F6").Value = 3
$ws2.Range("F7").Value = 824
$ws2.Range("F9").Value = 33

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 705
$ws4.Range("F5").Value = 11323
$ws4.Range("F6").Value = 1759
$ws4.Range("F7").Value = 84
$ws4.Range("F8").Value = 489
$ws4.Range("F9").Value = 781
$ws4.Range("F10").Value = 2362
$ws4.Range("F11").Value = 712
$ws4.Range("F12").Value = 945
$ws4.Range("F14").Value = 546
$ws4.Range("F15").Value = 409
$ws4.Range("F16").Value = 454
$ws4.Range("F17").Value = 418
$ws4.Range("F19").Value = 3
$ws4.Range("F20").Value = 581
$ws4.Range("F21").Value = 41
$ws4.Range("F22").Value = 824
$ws4.Range("F23").Value = 934
$ws4.Range("F24").Value = 406
$ws4.Range("F25").Value = 599
$ws4.Range("F26").Value = 949
$ws4.Range("F27").Value = 180
$ws4.Range("F28").Value = 893
$ws4.Range("F30").Value = 181
$ws4.Range("F34").Value = 627
$ws4.Range("F35").Value = 128
$ws4.Range("F36").Value = 33
$ws4.Range("F37").Value = 71
$ws4.Range("F38").Value = 294
